$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "kapelusz"
$ws.Range("B2").Value = "hat"
$ws.Range("B2").Select() | Out-Null
